# serie_precios_pgm.xlsx -- sep-2022 update
# - cleans up the redundant "applyFont only" style that used to sit on the
#   header row + the first three data rows (it rendered identically to the
#   unstyled default, so this just drops the explicit style reference)
# - appends the Sept-2022 data point (row 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cleanup: drop the redundant no-op style (header row + B2:C4) ---
$ws.Range("A1:C1").ClearFormats()
$ws.Range("B2:C4").ClearFormats()

# --- new row: 2022 Septiembre ---
$ws.Range("A9").Value = 44805
$ws.Range("B9").Value = 25.78
$ws.Range("C9").Value = 27.08

# give the new date cell the same look as the rest of column A (date format)
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# leave the selection where it ended up after entering the new row
$ws.Range("F10").Select() | Out-Null
